# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 509
$ws1.Range("F6").Value = 2708
$ws1.Range("F7").Value = 189
$ws1.Range("F10").Value = 1573
$ws1.Range("F11").Value = 551
$ws1.Range("F13").Value = 341
$ws1.Range("F18").Value = 216
$ws1.Range("F21").Value = 14
$ws1.Range("F22").Value = 213
$ws1.Range("F24").Value = 1740
$ws1.Range("F27").Value = 72
$ws1.Range("F28").Value = 570
$ws1.Range("F31").Value = 444

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 509
$ws4.Range("F7").Value = 2708
$ws4.Range("F8").Value = 189
$ws4.Range("F11").Value = 1573
$ws4.Range("F12").Value = 551
$ws4.Range("F14").Value = 341
$ws4.Range("F19").Value = 216
$ws4.Range("F22").Value = 14
$ws4.Range("F23").Value = 213
$ws4.Range("F25").Value = 1740
$ws4.Range("F28").Value = 72
$ws4.Range("F29").Value = 570
$ws4.Range("F32").Value = 444
